# feat: add 2022-Q1 data
#
# 1. Insert a new "2022-Q1" worksheet right before the "总计" (Total) sheet,
#    matching the layout of the other quarterly fund-holding sheets (we
#    clone the "2021-Q4" sheet so sheet-level properties/styles/margins
#    come along for free, then overwrite its cell contents).
# 2. Prepend a new "2022-Q1" row to the "总计" summary sheet, shifting the
#    existing rows down and renumbering the index column.

# Helper: force a cell to be stored as TEXT even when its value looks
# numeric (e.g. fund codes with leading zeros, "4.06", "0.0698", ...),
# then restore the default "Normal" style so no stray number format is
# left attached to the cell (matches the plain, unstyled text cells in
# the target workbook).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$firstSheet = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Step 1: create the "2022-Q1" sheet, positioned right before "总计"
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")

$template.Copy($totalSheet)
$q1 = $wb.Worksheets.Item("2021-Q4 (2)")
$q1.Name = "2022-Q1"

# Header row (labels differ from the quarterly fund sheets' template)
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Index column (0-based row counter) already carries the right style
# from the cloned template; just refresh the values (no-op here, but
# kept explicit/defensive).
$q1.Range("A2").Value = 0
$q1.Range("A3").Value = 1
$q1.Range("A4").Value = 2
$q1.Range("A5").Value = 3
$q1.Range("A6").Value = 4

# Data rows. Columns B (fund code) and D,E,F,G (numeric-looking
# percentages/amounts) are stored as TEXT in the source data, so they
# go through Set-TextValue; C (fund name) is plain text; H (rank) is a
# genuine number.
Set-TextValue $q1.Range("B2") "004818"
$q1.Range("C2").Value = "国寿安保目标策略灵活配置混合A"
Set-TextValue $q1.Range("D2") "4.06"
Set-TextValue $q1.Range("E2") "36.45"
Set-TextValue $q1.Range("F2") "1.72"
Set-TextValue $q1.Range("G2") "0.0698"
$q1.Range("H2").Value = 9

Set-TextValue $q1.Range("B3") "000458"
$q1.Range("C3").Value = "英大领先回报混合"
Set-TextValue $q1.Range("D3") "1.11"
Set-TextValue $q1.Range("E3") "93.55"
Set-TextValue $q1.Range("F3") "2.29"
Set-TextValue $q1.Range("G3") "0.0254"
$q1.Range("H3").Value = 6

Set-TextValue $q1.Range("B4") "001270"
$q1.Range("C4").Value = "英大灵活配置混合A"
Set-TextValue $q1.Range("D4") "0.56"
Set-TextValue $q1.Range("E4") "93.18"
Set-TextValue $q1.Range("F4") "1.98"
Set-TextValue $q1.Range("G4") "0.0111"
$q1.Range("H4").Value = 7

Set-TextValue $q1.Range("B5") "004819"
$q1.Range("C5").Value = "国寿安保目标策略灵活配置混合C"
Set-TextValue $q1.Range("D5") "0.57"
Set-TextValue $q1.Range("E5") "36.45"
Set-TextValue $q1.Range("F5") "1.72"
Set-TextValue $q1.Range("G5") "0.0098"
$q1.Range("H5").Value = 9

Set-TextValue $q1.Range("B6") "001271"
$q1.Range("C6").Value = "英大灵活配置混合B"
Set-TextValue $q1.Range("D6") "0.28"
Set-TextValue $q1.Range("E6") "93.18"
Set-TextValue $q1.Range("F6") "1.98"
Set-TextValue $q1.Range("G6") "0.0055"
$q1.Range("H6").Value = 7

# ---------------------------------------------------------------------
# Step 2: prepend a "2022-Q1" row to the "总计" sheet, shifting the rest
#          of the rows down by one and renumbering column A.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# The new row 7 needs the same "index column" style (s="2") that rows
# A2:A6 already carry; copy it down from A6 before the row gets a value.
$total.Range("A6").Copy()
$total.Range("A7").PasteSpecial(-4122)

for ($r = 6; $r -ge 2; $r--) {
    $srcRow = $r - 1
    $total.Cells.Item($r, 2).Value = $total.Cells.Item($srcRow, 2).Value2
    $total.Cells.Item($r, 3).Value = $total.Cells.Item($srcRow, 3).Value2
    $total.Cells.Item($r, 4).Value = $total.Cells.Item($srcRow, 4).Value2
}

$total.Range("B7").Value = "2020-Q4"
$total.Range("C7").Value = 5
$total.Range("D7").Value = 0.07000000000000001

$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 5
$total.Range("D2").Value = 0.12

for ($r = 2; $r -le 7; $r++) {
    $total.Cells.Item($r, 1).Value = $r - 2
}

# Restore the original active sheet/tab (creating/renaming sheets above
# shifts the active tab to whichever sheet was touched last).
$firstSheet.Activate()
